$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '45.213.05'
$ws.Cells.Item(2, 5).Value = '  -3.16%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.380.04'
$ws.Cells.Item(3, 5).Value = '  +3.53%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.10%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '296.98'
$ws.Cells.Item(5, 5).Value = '  -2.41%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '95.15'
$ws.Cells.Item(6, 5).Value = '  -6.28%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.558'
$ws.Cells.Item(7, 5).Value = '  -1.80%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '1.00'
$ws.Cells.Item(8, 5).Value = '  +0.01%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.504'
$ws.Cells.Item(9, 5).Value = '  -3.35%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '34.37'
$ws.Cells.Item(10, 5).Value = '  -5.78%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0779'
$ws.Cells.Item(11, 5).Value = '  -1.36%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '6.99'
$ws.Cells.Item(12, 5).Value = '  -4.66%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +0.59%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '2.761.01'
$ws.Cells.Item(14, 5).Value = '  +4.18%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '2.392.67'
$ws.Cells.Item(15, 5).Value = '  +3.97%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '14.00'
$ws.Cells.Item(16, 5).Value = '  +1.16%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.825'
$ws.Cells.Item(17, 5).Value = '  +1.79%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '45.177.35'
$ws.Cells.Item(18, 5).Value = '  -3.21%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '12.45'
$ws.Cells.Item(19, 5).Value = '  -3.97%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '0.0₃0934'
$ws.Cells.Item(20, 5).Value = '  -0.88%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '6.10'
$ws.Cells.Item(21, 5).Value = '  +1.69%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '66.52'
$ws.Cells.Item(22, 5).Value = '  +0.99%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '239.04'
$ws.Cells.Item(23, 5).Value = '  -4.12%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.78'
$ws.Cells.Item(24, 5).Value = '  -3.77%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.00%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '1.89'
$ws.Cells.Item(26, 5).Value = '  -1.68%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +0.52%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '37.71'
$ws.Cells.Item(28, 5).Value = '  -10.75%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '9.60'
$ws.Cells.Item(29, 5).Value = '  -2.88%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '3.79'
$ws.Cells.Item(30, 5).Value = '  +15.58%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '20.96'
$ws.Cells.Item(31, 5).Value = '  +4.76%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '149.17'
$ws.Cells.Item(32, 5).Value = '  +1.12%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -3.64%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '5.44'
$ws.Cells.Item(34, 5).Value = '  -2.76%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.0761'
$ws.Cells.Item(35, 5).Value = '  -3.58%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'Kaspa'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.112'
$ws.Cells.Item(36, 5).Value = '  -3.06%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'ARBITRUM'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '1.94'
$ws.Cells.Item(37, 5).Value = '  +9.94%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.115'
$ws.Cells.Item(38, 5).Value = '  -2.01%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '14.85'
$ws.Cells.Item(39, 5).Value = '  -6.96%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.76'
$ws.Cells.Item(40, 5).Value = '  -5.27%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.0295'
$ws.Cells.Item(41, 5).Value = '  -1.75%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.18'
$ws.Cells.Item(42, 5).Value = '  -4.88%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '1.944.80'
$ws.Cells.Item(43, 5).Value = '  +6.61%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +0.04%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '90.08'
$ws.Cells.Item(45, 5).Value = '  +1.56%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '1.70'
$ws.Cells.Item(46, 5).Value = '  -12.81%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'FraxShare'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '8.64'
$ws.Cells.Item(47, 5).Value = '  +9.78%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '15.48'
$ws.Cells.Item(48, 5).Value = '  +16.91%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '100.49'
$ws.Cells.Item(49, 5).Value = '  +5.26%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(50, 4).Value = '2.624.74'
$ws.Cells.Item(50, 5).Value = '  +3.95%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'Algorand'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.184'
$ws.Cells.Item(51, 5).Value = '  -5.11%  '
